# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 23:39"

# --- Country ranking swaps (rows whose country overtook its neighbour) ---
# Costa Rica overtakes Japon
$ws.Range("A49").Value = "Costa Rica"
$ws.Range("A50").Value = "Japon"

# Trinidad y Tobago overtakes Congo
$ws.Range("A129").Value = "Trinidad yTobago"
$ws.Range("A130").Value = "Congo"

# --- Updated case figures ---
$ws.Range("B4").Value = 8081232
$ws.Range("C4").Value = 42906
$ws.Range("D4").Value = 5213465
$ws.Range("E4").Value = 2647074
$ws.Range("G4").Value = 675
$ws.Range("H4").Value = 220693

$ws.Range("B27").Value = 296652
$ws.Range("C27").Value = 2621
$ws.Range("D27").Value = 245811
$ws.Range("E27").Value = 48786
$ws.Range("G27").Value = 34
$ws.Range("H27").Value = 2055

$ws.Range("B49").Value = 90238
$ws.Range("C49").Value = 1015
$ws.Range("D49").Value = 53670
$ws.Range("E49").Value = 35444
$ws.Range("G49").Value = 16
$ws.Range("H49").Value = 1124

$ws.Range("B50").Value = 89673
$ws.Range("C50").Value = 326
$ws.Range("D50").Value = 82621
$ws.Range("E50").Value = 5418
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 1634

$ws.Range("B66").Value = 53399
$ws.Range("C66").Value = 174
$ws.Range("D66").Value = 37492
$ws.Range("E66").Value = 14089
$ws.Range("G66").Value = 9
$ws.Range("H66").Value = 1818

$ws.Range("B129").Value = 5127
$ws.Range("C129").Value = 11
$ws.Range("D129").Value = 3367
$ws.Range("E129").Value = 1667
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 93

$ws.Range("B130").Value = 5118
$ws.Range("D130").Value = 3887
$ws.Range("E130").Value = 1141
$ws.Range("H130").Value = 90

$ws.Range("B146").Value = 3565
$ws.Range("C146").Value = 44
$ws.Range("D146").Value = 2435
$ws.Range("E146").Value = 1024
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 106

$ws.Range("B148").Value = 3297
$ws.Range("C148").Value = 1
$ws.Range("D148").Value = 2544
$ws.Range("E148").Value = 621

$ws.Range("B159").Value = 2130
$ws.Range("C159").Value = 83
$ws.Range("E159").Value = 661

$ws.Range("B160").Value = 2053
$ws.Range("C160").Value = 1
$ws.Range("E160").Value = 128

$ws.Range("B167").Value = 1203
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 10

$ws.Range("B184").Value = 407
$ws.Range("C184").Value = 3
$ws.Range("D184").Value = 361
$ws.Range("E184").Value = 36
